# Fixed #366 User content is lost after two generation without edition.
#
# Turns the two m2doc "user content" markers (m:usercontent / m:endusercontent),
# which are stored as simple fields (<w:fldSimple w:instr="..."/>), into
# complex fields (begin / instrText / separate / end run sequence) so that a
# second generation does not lose the user-edited content between them.

$d = $word.ActiveDocument

function Get-ContainingParagraph {
    param($doc, $range)

    $target = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Start -le $range.Start -and $range.End -le $p.Range.End) {
            $target = $p
        }
    }
    return $target
}

function Convert-SimpleFieldToComplex {
    param($doc, [string]$instrText)

    $targetField = $null
    foreach ($f in $doc.Fields) {
        if ($f.Code.Text.Trim() -eq $instrText.Trim()) {
            $targetField = $f
        }
    }
    if ($targetField -eq $null) {
        throw "Convert-SimpleFieldToComplex: field not found for '$instrText'"
    }

    $code = $targetField.Code
    $instr = $code.Text
    $para = Get-ContainingParagraph $doc $code

    # Pull the paragraph's own OOXML so we keep its attributes (rsids, ...)
    # and any sibling content (bookmarks, ...) untouched.
    $paraXml = $para.Range.WordOpenXML
    $paraMatch = [regex]::Match(
        $paraXml,
        "<w:p\b([^>]*)>(.*?)</w:p>",
        [System.Text.RegularExpressions.RegexOptions]::Singleline)

    $paraAttrs = $paraMatch.Groups[1].Value
    # WordOpenXML() stamps paragraph identity attributes that do not exist
    # in the source document; drop them so we do not introduce them.
    $paraAttrs = [regex]::Replace($paraAttrs, '\s*w14:paraId="[^"]*"', '')
    $paraAttrs = [regex]::Replace($paraAttrs, '\s*w14:textId="[^"]*"', '')

    $innerContent = $paraMatch.Groups[2].Value
    # Drop the old <w:fldSimple .../> (self closed or with content), keep
    # everything else in the paragraph (e.g. bookmarkStart/bookmarkEnd).
    $preserved = [regex]::Replace(
        $innerContent,
        "<w:fldSimple\b[^>]*(/>|>.*?</w:fldSimple>)",
        "",
        [System.Text.RegularExpressions.RegexOptions]::Singleline)

    $instrEsc = $instr -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'

    $complexField = ('<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText>' + $instrEsc + '</w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>')

    $newParaInner = $preserved + $complexField

    $xml = ('<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p' + $paraAttrs + '>' + $newParaInner + '</w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>')

    $targetRange = $doc.Range($code.Start, $code.End)
    $targetRange.InsertXML($xml)
}

Convert-SimpleFieldToComplex $d "m:usercontent  zone1"
Convert-SimpleFieldToComplex $d "m:endusercontent"
